# ===========================================================================
# Add files via upload
#
# The underlying experiment was re-run and produced a few more recorded
# apple-eating episodes. This adds those new result rows to the BFS, DFS
# and Simulated_Annealing sheets (each sheet lists, per apple eaten: the
# apple index, its (row, col) position, whether it was reached, the number
# of states explored, and the solution path), and then refreshes the
# "Average explored states:" SUM/COUNTIF formula on row 15 of every sheet
# so its range covers all of the data rows that are actually present.
# ===========================================================================

$xlCenter = -4108

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# BFS: apple runs #5-#10 recorded as new rows 7-12
# ---------------------------------------------------------------------------
$wsBFS = $wb.Worksheets.Item("BFS")


    # Row 7
    $wsBFS.Cells.Item(7, 1).Value = 5
    $wsBFS.Cells.Item(7, 2).Value = '((4, 13))'
    $wsBFS.Cells.Item(7, 3).Value = $true
    $wsBFS.Cells.Item(7, 4).Value = 275
    $wsBFS.Cells.Item(7, 5).Value = '(22, 7) -> (21, 7) -> (21, 8) -> (21, 9) -> (21, 10) -> (20, 10) -> (19, 10) -> (18, 10) -> (17, 10) -> (16, 10) -> (15, 10) -> (14, 10) -> (13, 10) -> (12, 10) -> (11, 10) -> (10, 10) -> (9, 10) -> (8, 10) -> (7, 10) -> (6, 10) -> (5, 10) -> (4, 10) -> (4, 11) -> (4, 12) -> (4, 13)'
    $wsBFS.Range("A7:D7").HorizontalAlignment = $xlCenter
    $wsBFS.Range("A7:E7").VerticalAlignment = $xlCenter
    $wsBFS.Range("E7").WrapText = $true

    # Row 8
    $wsBFS.Cells.Item(8, 1).Value = 6
    $wsBFS.Cells.Item(8, 2).Value = '((19, 1))'
    $wsBFS.Cells.Item(8, 3).Value = $true
    $wsBFS.Cells.Item(8, 4).Value = 309
    $wsBFS.Cells.Item(8, 5).Value = '(4, 13) -> (5, 13) -> (6, 13) -> (7, 13) -> (8, 13) -> (9, 13) -> (10, 13) -> (11, 13) -> (12, 13) -> (13, 13) -> (14, 13) -> (15, 13) -> (16, 13) -> (17, 13) -> (18, 13) -> (19, 13) -> (19, 12) -> (19, 11) -> (19, 10) -> (19, 9) -> (19, 8) -> (19, 7) -> (19, 6) -> (19, 5) -> (19, 4) -> (19, 3) -> (19, 2) -> (19, 1)'
    $wsBFS.Range("A8:D8").HorizontalAlignment = $xlCenter
    $wsBFS.Range("A8:E8").VerticalAlignment = $xlCenter
    $wsBFS.Range("E8").WrapText = $true

    # Row 9
    $wsBFS.Cells.Item(9, 1).Value = 7
    $wsBFS.Cells.Item(9, 2).Value = '((5, 5))'
    $wsBFS.Cells.Item(9, 3).Value = $true
    $wsBFS.Cells.Item(9, 4).Value = 212
    $wsBFS.Cells.Item(9, 5).Value = '(19, 1) -> (18, 1) -> (17, 1) -> (16, 1) -> (15, 1) -> (14, 1) -> (13, 1) -> (12, 1) -> (11, 1) -> (10, 1) -> (10, 2) -> (10, 3) -> (9, 3) -> (8, 3) -> (7, 3) -> (6, 3) -> (5, 3) -> (5, 4) -> (5, 5)'
    $wsBFS.Range("A9:D9").HorizontalAlignment = $xlCenter
    $wsBFS.Range("A9:E9").VerticalAlignment = $xlCenter
    $wsBFS.Range("E9").WrapText = $true

    # Row 10
    $wsBFS.Cells.Item(10, 1).Value = 8
    $wsBFS.Cells.Item(10, 2).Value = '((12, 14))'
    $wsBFS.Cells.Item(10, 3).Value = $true
    $wsBFS.Cells.Item(10, 4).Value = 239
    $wsBFS.Cells.Item(10, 5).Value = '(5, 5) -> (6, 5) -> (7, 5) -> (8, 5) -> (9, 5) -> (10, 5) -> (11, 5) -> (12, 5) -> (12, 6) -> (12, 7) -> (12, 8) -> (12, 9) -> (12, 10) -> (12, 11) -> (12, 12) -> (12, 13) -> (12, 14)'
    $wsBFS.Range("A10:D10").HorizontalAlignment = $xlCenter
    $wsBFS.Range("A10:E10").VerticalAlignment = $xlCenter
    $wsBFS.Range("E10").WrapText = $true

    # Row 11
    $wsBFS.Cells.Item(11, 1).Value = 9
    $wsBFS.Cells.Item(11, 2).Value = '((14, 8))'
    $wsBFS.Cells.Item(11, 3).Value = $true
    $wsBFS.Cells.Item(11, 4).Value = 56
    $wsBFS.Cells.Item(11, 5).Value = '(12, 14) -> (13, 14) -> (13, 13) -> (13, 12) -> (13, 11) -> (14, 11) -> (14, 10) -> (14, 9) -> (14, 8)'
    $wsBFS.Range("A11:D11").HorizontalAlignment = $xlCenter
    $wsBFS.Range("A11:E11").VerticalAlignment = $xlCenter
    $wsBFS.Range("E11").WrapText = $true

    # Row 12
    $wsBFS.Cells.Item(12, 1).Value = 10
    $wsBFS.Cells.Item(12, 2).Value = '((22, 13))'
    $wsBFS.Cells.Item(12, 3).Value = $true
    $wsBFS.Cells.Item(12, 4).Value = 203
    $wsBFS.Cells.Item(12, 5).Value = '(14, 8) -> (15, 8) -> (16, 8) -> (17, 8) -> (18, 8) -> (19, 8) -> (19, 9) -> (19, 10) -> (20, 10) -> (21, 10) -> (22, 10) -> (22, 11) -> (22, 12) -> (22, 13)'
    $wsBFS.Range("A12:D12").HorizontalAlignment = $xlCenter
    $wsBFS.Range("A12:E12").VerticalAlignment = $xlCenter
    $wsBFS.Range("E12").WrapText = $true

$wsBFS.Range("D15").Formula = "=SUM(D3:D12) / COUNTIF(C3:C12, True)"


# ---------------------------------------------------------------------------
# DFS: apple runs #2-#10 recorded as new rows 4-12
# ---------------------------------------------------------------------------
$wsDFS = $wb.Worksheets.Item("DFS")


    # Row 4
    $wsDFS.Cells.Item(4, 1).Value = 2
    $wsDFS.Cells.Item(4, 2).Value = '((1, 1))'
    $wsDFS.Cells.Item(4, 3).Value = $true
    $wsDFS.Cells.Item(4, 4).Value = 79
    $wsDFS.Cells.Item(4, 5).Value = '(12, 7) -> (12, 8) -> (12, 9) -> (12, 10) -> (12, 11) -> (12, 12) -> (12, 13) -> (12, 14) -> (11, 14) -> (10, 14) -> (10, 13) -> (10, 12) -> (10, 11) -> (10, 10) -> (10, 9) -> (10, 8) -> (10, 7) -> (10, 6) -> (10, 5) -> (10, 4) -> (10, 3) -> (10, 2) -> (10, 1) -> (10, 0) -> (9, 0) -> (8, 0) -> (8, 1) -> (7, 1) -> (6, 1) -> (6, 0) -> (5, 0) -> (4, 0) -> (4, 1) -> (4, 2) -> (4, 3) -> (4, 4) -> (4, 5) -> (4, 6) -> (4, 7) -> (4, 8) -> (4, 9) -> (4, 10) -> (4, 11) -> (4, 12) -> (4, 13) -> (4, 14) -> (3, 14) -> (2, 14) -> (2, 13) -> (2, 12) -> (2, 11) -> (2, 10) -> (2, 9) -> (2, 8) -> (2, 7) -> (2, 6) -> (2, 5) -> (2, 4) -> (2, 3) -> (2, 2) -> (2, 1) -> (1, 1)'
    $wsDFS.Range("A4:D4").HorizontalAlignment = $xlCenter
    $wsDFS.Range("A4:E4").VerticalAlignment = $xlCenter
    $wsDFS.Range("E4").WrapText = $true

    # Row 5
    $wsDFS.Cells.Item(5, 1).Value = 3
    $wsDFS.Cells.Item(5, 2).Value = '((11, 3))'
    $wsDFS.Cells.Item(5, 3).Value = $true
    $wsDFS.Cells.Item(5, 4).Value = 275
    $wsDFS.Cells.Item(5, 5).Value = '(1, 1) -> (1, 0) -> (2, 0) -> (3, 0) -> (3, 1) -> (3, 2) -> (3, 3) -> (3, 4) -> (2, 4) -> (2, 5) -> (2, 6) -> (2, 7) -> (2, 8) -> (2, 9) -> (2, 10) -> (2, 11) -> (2, 12) -> (2, 13) -> (2, 14) -> (3, 14) -> (4, 14) -> (4, 13) -> (4, 12) -> (4, 11) -> (4, 10) -> (4, 9) -> (4, 8) -> (4, 7) -> (4, 6) -> (4, 5) -> (5, 5) -> (5, 4) -> (5, 3) -> (6, 3) -> (7, 3) -> (7, 4) -> (7, 5) -> (7, 6) -> (7, 7) -> (7, 8) -> (7, 9) -> (7, 10) -> (7, 11) -> (8, 11) -> (9, 11) -> (9, 10) -> (9, 9) -> (9, 8) -> (9, 7) -> (9, 6) -> (9, 5) -> (9, 4) -> (9, 3) -> (10, 3) -> (11, 3)'
    $wsDFS.Range("A5:D5").HorizontalAlignment = $xlCenter
    $wsDFS.Range("A5:E5").VerticalAlignment = $xlCenter
    $wsDFS.Range("E5").WrapText = $true

    # Row 6
    $wsDFS.Cells.Item(6, 1).Value = 4
    $wsDFS.Cells.Item(6, 2).Value = '((22, 7))'
    $wsDFS.Cells.Item(6, 3).Value = $true
    $wsDFS.Cells.Item(6, 4).Value = 249
    $wsDFS.Cells.Item(6, 5).Value = '(11, 3) -> (11, 2) -> (11, 1) -> (11, 0) -> (10, 0) -> (9, 0) -> (9, 1) -> (8, 1) -> (7, 1) -> (7, 0) -> (6, 0) -> (5, 0) -> (5, 1) -> (4, 1) -> (4, 2) -> (4, 3) -> (4, 4) -> (4, 5) -> (4, 6) -> (4, 7) -> (4, 8) -> (4, 9) -> (4, 10) -> (4, 11) -> (4, 12) -> (4, 13) -> (4, 14) -> (5, 14) -> (6, 14) -> (6, 13) -> (7, 13) -> (8, 13) -> (8, 14) -> (9, 14) -> (10, 14) -> (10, 13) -> (10, 12) -> (10, 11) -> (10, 10) -> (10, 9) -> (10, 8) -> (10, 7) -> (10, 6) -> (10, 5) -> (11, 5) -> (12, 5) -> (12, 4) -> (13, 4) -> (13, 3) -> (13, 2) -> (13, 1) -> (13, 0) -> (14, 0) -> (15, 0) -> (15, 1) -> (16, 1) -> (17, 1) -> (17, 0) -> (18, 0) -> (19, 0) -> (19, 1) -> (19, 2) -> (19, 3) -> (19, 4) -> (19, 5) -> (19, 6) -> (19, 7) -> (19, 8) -> (19, 9) -> (19, 10) -> (19, 11) -> (19, 12) -> (19, 13) -> (19, 14) -> (20, 14) -> (21, 14) -> (21, 13) -> (21, 12) -> (21, 11) -> (21, 10) -> (21, 9) -> (21, 8) -> (21, 7) -> (22, 7)'
    $wsDFS.Range("A6:D6").HorizontalAlignment = $xlCenter
    $wsDFS.Range("A6:E6").VerticalAlignment = $xlCenter
    $wsDFS.Range("E6").WrapText = $true

    # Row 7
    $wsDFS.Cells.Item(7, 1).Value = 5
    $wsDFS.Cells.Item(7, 2).Value = '((4, 13))'
    $wsDFS.Cells.Item(7, 3).Value = $true
    $wsDFS.Cells.Item(7, 4).Value = 73
    $wsDFS.Cells.Item(7, 5).Value = '(22, 7) -> (22, 6) -> (22, 5) -> (22, 4) -> (22, 3) -> (22, 2) -> (22, 1) -> (22, 0) -> (21, 0) -> (20, 0) -> (20, 1) -> (20, 2) -> (20, 3) -> (20, 4) -> (19, 4) -> (19, 5) -> (19, 6) -> (19, 7) -> (19, 8) -> (19, 9) -> (19, 10) -> (19, 11) -> (19, 12) -> (19, 13) -> (19, 14) -> (18, 14) -> (17, 14) -> (17, 13) -> (16, 13) -> (15, 13) -> (15, 14) -> (14, 14) -> (13, 14) -> (13, 13) -> (13, 12) -> (13, 11) -> (13, 10) -> (13, 9) -> (13, 8) -> (13, 7) -> (13, 6) -> (13, 5) -> (13, 4) -> (13, 3) -> (13, 2) -> (13, 1) -> (13, 0) -> (12, 0) -> (11, 0) -> (11, 1) -> (11, 2) -> (11, 3) -> (11, 4) -> (11, 5) -> (11, 6) -> (11, 7) -> (11, 8) -> (11, 9) -> (11, 10) -> (11, 11) -> (11, 12) -> (11, 13) -> (11, 14) -> (10, 14) -> (9, 14) -> (9, 13) -> (8, 13) -> (7, 13) -> (7, 14) -> (6, 14) -> (5, 14) -> (5, 13) -> (4, 13)'
    $wsDFS.Range("A7:D7").HorizontalAlignment = $xlCenter
    $wsDFS.Range("A7:E7").VerticalAlignment = $xlCenter
    $wsDFS.Range("E7").WrapText = $true

    # Row 8
    $wsDFS.Cells.Item(8, 1).Value = 6
    $wsDFS.Cells.Item(8, 2).Value = '((19, 1))'
    $wsDFS.Cells.Item(8, 3).Value = $true
    $wsDFS.Cells.Item(8, 4).Value = 112
    $wsDFS.Cells.Item(8, 5).Value = '(4, 13) -> (4, 12) -> (4, 11) -> (4, 10) -> (4, 9) -> (4, 8) -> (4, 7) -> (4, 6) -> (4, 5) -> (4, 4) -> (4, 3) -> (4, 2) -> (4, 1) -> (4, 0) -> (5, 0) -> (6, 0) -> (6, 1) -> (7, 1) -> (8, 1) -> (8, 0) -> (9, 0) -> (10, 0) -> (10, 1) -> (10, 2) -> (10, 3) -> (10, 4) -> (10, 5) -> (10, 6) -> (10, 7) -> (10, 8) -> (10, 9) -> (10, 10) -> (10, 11) -> (10, 12) -> (10, 13) -> (10, 14) -> (11, 14) -> (12, 14) -> (12, 13) -> (12, 12) -> (12, 11) -> (12, 10) -> (12, 9) -> (12, 8) -> (12, 7) -> (12, 6) -> (12, 5) -> (12, 4) -> (12, 3) -> (12, 2) -> (12, 1) -> (12, 0) -> (13, 0) -> (14, 0) -> (14, 1) -> (15, 1) -> (16, 1) -> (16, 0) -> (17, 0) -> (18, 0) -> (18, 1) -> (19, 1)'
    $wsDFS.Range("A8:D8").HorizontalAlignment = $xlCenter
    $wsDFS.Range("A8:E8").VerticalAlignment = $xlCenter
    $wsDFS.Range("E8").WrapText = $true

    # Row 9
    $wsDFS.Cells.Item(9, 1).Value = 7
    $wsDFS.Cells.Item(9, 2).Value = '((5, 5))'
    $wsDFS.Cells.Item(9, 3).Value = $true
    $wsDFS.Cells.Item(9, 4).Value = 195
    $wsDFS.Cells.Item(9, 5).Value = '(19, 1) -> (19, 0) -> (20, 0) -> (21, 0) -> (21, 1) -> (21, 2) -> (21, 3) -> (21, 4) -> (21, 5) -> (21, 6) -> (21, 7) -> (21, 8) -> (21, 9) -> (21, 10) -> (21, 11) -> (21, 12) -> (21, 13) -> (21, 14) -> (20, 14) -> (19, 14) -> (19, 13) -> (19, 12) -> (19, 11) -> (19, 10) -> (19, 9) -> (19, 8) -> (19, 7) -> (19, 6) -> (19, 5) -> (19, 4) -> (19, 3) -> (18, 3) -> (17, 3) -> (17, 4) -> (17, 5) -> (17, 6) -> (17, 7) -> (17, 8) -> (17, 9) -> (17, 10) -> (17, 11) -> (16, 11) -> (15, 11) -> (15, 10) -> (15, 9) -> (15, 8) -> (15, 7) -> (15, 6) -> (15, 5) -> (15, 4) -> (15, 3) -> (14, 3) -> (13, 3) -> (13, 2) -> (13, 1) -> (13, 0) -> (12, 0) -> (11, 0) -> (11, 1) -> (11, 2) -> (11, 3) -> (11, 4) -> (11, 5) -> (11, 6) -> (11, 7) -> (11, 8) -> (11, 9) -> (11, 10) -> (11, 11) -> (11, 12) -> (11, 13) -> (11, 14) -> (10, 14) -> (9, 14) -> (9, 13) -> (8, 13) -> (7, 13) -> (7, 14) -> (6, 14) -> (5, 14) -> (5, 13) -> (4, 13) -> (4, 12) -> (4, 11) -> (4, 10) -> (4, 9) -> (4, 8) -> (4, 7) -> (4, 6) -> (4, 5) -> (5, 5)'
    $wsDFS.Range("A9:D9").HorizontalAlignment = $xlCenter
    $wsDFS.Range("A9:E9").VerticalAlignment = $xlCenter
    $wsDFS.Range("E9").WrapText = $true

    # Row 10
    $wsDFS.Cells.Item(10, 1).Value = 8
    $wsDFS.Cells.Item(10, 2).Value = '((12, 14))'
    $wsDFS.Cells.Item(10, 3).Value = $true
    $wsDFS.Cells.Item(10, 4).Value = 87
    $wsDFS.Cells.Item(10, 5).Value = '(5, 5) -> (5, 4) -> (5, 3) -> (4, 3) -> (4, 2) -> (4, 1) -> (4, 0) -> (3, 0) -> (2, 0) -> (2, 1) -> (2, 2) -> (2, 3) -> (2, 4) -> (2, 5) -> (2, 6) -> (2, 7) -> (2, 8) -> (2, 9) -> (2, 10) -> (2, 11) -> (2, 12) -> (2, 13) -> (2, 14) -> (3, 14) -> (4, 14) -> (4, 13) -> (5, 13) -> (6, 13) -> (6, 14) -> (7, 14) -> (8, 14) -> (8, 13) -> (9, 13) -> (10, 13) -> (10, 12) -> (10, 11) -> (10, 10) -> (10, 9) -> (10, 8) -> (10, 7) -> (10, 6) -> (10, 5) -> (10, 4) -> (10, 3) -> (10, 2) -> (10, 1) -> (10, 0) -> (11, 0) -> (12, 0) -> (12, 1) -> (12, 2) -> (12, 3) -> (12, 4) -> (12, 5) -> (12, 6) -> (12, 7) -> (12, 8) -> (12, 9) -> (12, 10) -> (12, 11) -> (12, 12) -> (12, 13) -> (12, 14)'
    $wsDFS.Range("A10:D10").HorizontalAlignment = $xlCenter
    $wsDFS.Range("A10:E10").VerticalAlignment = $xlCenter
    $wsDFS.Range("E10").WrapText = $true

    # Row 11
    $wsDFS.Cells.Item(11, 1).Value = 9
    $wsDFS.Cells.Item(11, 2).Value = '((14, 8))'
    $wsDFS.Cells.Item(11, 3).Value = $true
    $wsDFS.Cells.Item(11, 4).Value = 295
    $wsDFS.Cells.Item(11, 5).Value = '(12, 14) -> (11, 14) -> (11, 13) -> (11, 12) -> (11, 11) -> (11, 10) -> (11, 9) -> (11, 8) -> (11, 7) -> (11, 6) -> (11, 5) -> (11, 4) -> (11, 3) -> (11, 2) -> (11, 1) -> (11, 0) -> (12, 0) -> (13, 0) -> (13, 1) -> (13, 2) -> (13, 3) -> (13, 4) -> (13, 5) -> (13, 6) -> (13, 7) -> (13, 8) -> (14, 8)'
    $wsDFS.Range("A11:D11").HorizontalAlignment = $xlCenter
    $wsDFS.Range("A11:E11").VerticalAlignment = $xlCenter
    $wsDFS.Range("E11").WrapText = $true

    # Row 12
    $wsDFS.Cells.Item(12, 1).Value = 10
    $wsDFS.Cells.Item(12, 2).Value = '((22, 13))'
    $wsDFS.Cells.Item(12, 3).Value = $true
    $wsDFS.Cells.Item(12, 4).Value = 227
    $wsDFS.Cells.Item(12, 5).Value = '(14, 8) -> (14, 7) -> (14, 6) -> (14, 5) -> (14, 4) -> (14, 3) -> (15, 3) -> (16, 3) -> (16, 4) -> (16, 5) -> (16, 6) -> (16, 7) -> (16, 8) -> (16, 9) -> (16, 10) -> (16, 11) -> (15, 11) -> (14, 11) -> (14, 10) -> (13, 10) -> (13, 9) -> (12, 9) -> (12, 8) -> (12, 7) -> (12, 6) -> (12, 5) -> (12, 4) -> (12, 3) -> (12, 2) -> (12, 1) -> (12, 0) -> (11, 0) -> (10, 0) -> (10, 1) -> (10, 2) -> (10, 3) -> (10, 4) -> (10, 5) -> (10, 6) -> (10, 7) -> (10, 8) -> (10, 9) -> (10, 10) -> (10, 11) -> (10, 12) -> (10, 13) -> (10, 14) -> (11, 14) -> (12, 14) -> (12, 13) -> (13, 13) -> (14, 13) -> (14, 14) -> (15, 14) -> (16, 14) -> (16, 13) -> (17, 13) -> (18, 13) -> (18, 14) -> (19, 14) -> (20, 14) -> (20, 13) -> (20, 12) -> (20, 11) -> (20, 10) -> (19, 10) -> (19, 9) -> (19, 8) -> (19, 7) -> (19, 6) -> (19, 5) -> (19, 4) -> (19, 3) -> (19, 2) -> (19, 1) -> (19, 0) -> (20, 0) -> (21, 0) -> (21, 1) -> (21, 2) -> (21, 3) -> (21, 4) -> (21, 5) -> (21, 6) -> (21, 7) -> (21, 8) -> (21, 9) -> (22, 9) -> (22, 10) -> (22, 11) -> (22, 12) -> (22, 13)'
    $wsDFS.Range("A12:D12").HorizontalAlignment = $xlCenter
    $wsDFS.Range("A12:E12").VerticalAlignment = $xlCenter
    $wsDFS.Range("E12").WrapText = $true

$wsDFS.Range("D15").Formula = "=SUM(D3:D12) / COUNTIF(C3:C12, True)"


# ---------------------------------------------------------------------------
# UCS, Greedy, Beam, Partially_Observable, Backtracking: the data rows
# already run down to row 12, but row 15's summary formula was left over
# from an earlier, shorter run (only covering D3:D11 / C3:C11) - extend it
# to include row 12 as well.
# ---------------------------------------------------------------------------
foreach ($name in @("UCS", "Greedy", "Beam", "Partially_Observable", "Backtracking")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("D15").Formula = "=SUM(D3:D12) / COUNTIF(C3:C12, True)"
}


# ---------------------------------------------------------------------------
# Simulated_Annealing: apple run #8 now has an extra intermediate hop
# recorded, so the old row 11 is replaced by two rows (new rows 11 & 12),
# and the former rows 12/13 shift down to rows 13/14.
# ---------------------------------------------------------------------------
$wsSA = $wb.Worksheets.Item("Simulated_Annealing")


    # Row 10
    $wsSA.Cells.Item(10, 1).Value = 8
    $wsSA.Cells.Item(10, 2).Value = '((12, 14))'
    $wsSA.Cells.Item(10, 3).Value = $false
    $wsSA.Cells.Item(10, 4).Value = 66
    $wsSA.Cells.Item(10, 5).Value = '(5, 5) -> (4, 5)'
    $wsSA.Range("A10:D10").HorizontalAlignment = $xlCenter
    $wsSA.Range("A10:E10").VerticalAlignment = $xlCenter
    $wsSA.Range("E10").WrapText = $true

    # Row 11
    $wsSA.Cells.Item(11, 1).Value = 8
    $wsSA.Cells.Item(11, 2).Value = '((12, 14))'
    $wsSA.Cells.Item(11, 3).Value = $false
    $wsSA.Cells.Item(11, 4).Value = 66
    $wsSA.Cells.Item(11, 5).Value = '(4, 5) -> (4, 4) -> (3, 4) -> (2, 4) -> (1, 4) -> (0, 4) -> (0, 3) -> (0, 2) -> (0, 1) -> (0, 0)'
    $wsSA.Range("A11:D11").HorizontalAlignment = $xlCenter
    $wsSA.Range("A11:E11").VerticalAlignment = $xlCenter
    $wsSA.Range("E11").WrapText = $true

    # Row 12
    $wsSA.Cells.Item(12, 1).Value = 8
    $wsSA.Cells.Item(12, 2).Value = '((12, 14))'
    $wsSA.Cells.Item(12, 3).Value = $true
    $wsSA.Cells.Item(12, 4).Value = 66
    $wsSA.Cells.Item(12, 5).Value = '(0, 0) -> (1, 0) -> (2, 0) -> (3, 0) -> (4, 0) -> (5, 0) -> (6, 0) -> (7, 0) -> (8, 0) -> (9, 0) -> (10, 0) -> (11, 0) -> (12, 0) -> (12, 1) -> (12, 2) -> (12, 3) -> (12, 4) -> (12, 5) -> (12, 6) -> (12, 7) -> (12, 8) -> (12, 9) -> (12, 10) -> (12, 11) -> (12, 12) -> (12, 13) -> (12, 14)'
    $wsSA.Range("A12:D12").HorizontalAlignment = $xlCenter
    $wsSA.Range("A12:E12").VerticalAlignment = $xlCenter
    $wsSA.Range("E12").WrapText = $true

    # Row 13
    $wsSA.Cells.Item(13, 1).Value = 9
    $wsSA.Cells.Item(13, 2).Value = '((14, 8))'
    $wsSA.Cells.Item(13, 3).Value = $true
    $wsSA.Cells.Item(13, 4).Value = 66
    $wsSA.Cells.Item(13, 5).Value = '(12, 14) -> (13, 14) -> (13, 13) -> (13, 12) -> (13, 11) -> (14, 11) -> (14, 10) -> (14, 9) -> (14, 8)'
    $wsSA.Range("A13:D13").HorizontalAlignment = $xlCenter
    $wsSA.Range("A13:E13").VerticalAlignment = $xlCenter
    $wsSA.Range("E13").WrapText = $true

    # Row 14
    $wsSA.Cells.Item(14, 1).Value = 10
    $wsSA.Cells.Item(14, 2).Value = '((22, 13))'
    $wsSA.Cells.Item(14, 3).Value = $true
    $wsSA.Cells.Item(14, 4).Value = 66
    $wsSA.Cells.Item(14, 5).Value = '(14, 8) -> (15, 8) -> (16, 8) -> (17, 8) -> (18, 8) -> (19, 8) -> (19, 9) -> (19, 10) -> (20, 10) -> (21, 10) -> (22, 10) -> (22, 11) -> (22, 12) -> (22, 13)'
    $wsSA.Range("A14:D14").HorizontalAlignment = $xlCenter
    $wsSA.Range("A14:E14").VerticalAlignment = $xlCenter
    $wsSA.Range("E14").WrapText = $true

$wsSA.Range("D15").Formula = "=SUM(D3:D14) / COUNTIF(C3:C14, True)"
